$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Rows 60 and 61 had their match data (columns F:V) swapped.
#    Columns A:E (index / pais / torneio / temporada / data_partida)
#    are identical between the two rows, so only F:V are touched.
#    (Use .Value2 for reads - .Value on this host only yields a
#    usable scalar when read directly into a Range.Value write,
#    not when captured into an intermediate variable.)
# ---------------------------------------------------------------

$swapCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
foreach ($c in $swapCols) {
    $v60 = $ws.Range("$c" + "60").Value2
    $v61 = $ws.Range("$c" + "61").Value2
    $ws.Range("$c" + "60").Value = $v61
    $ws.Range("$c" + "61").Value = $v60
}

# ---------------------------------------------------------------
# 2) Append 8 new match rows (130-137). First clone the row-129
#    formatting (bold/bordered/centered style on column A and the
#    datetime number format on column E) down onto the new rows,
#    then fill in the values.
# ---------------------------------------------------------------

$ws.Range("A129:V129").Copy()
$ws.Range("A130:V137").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @{r=130; idx=129; dataPartida=45282.90625; home="Os Belenenses"; homeGols=0; away="Santa Clara"; awayGols=0;
      hOpen=3.91; hOpenDt="17/12/2023 13:12"; hClose=5.64; hCloseDt="22/12/2023 21:37";
      dOpen=3.34; dOpenDt="17/12/2023 13:12"; dClose=3.67; dCloseDt="22/12/2023 21:37";
      aOpen=1.95; aOpenDt="17/12/2023 13:12"; aClose=1.71; aCloseDt="22/12/2023 21:43";
      url="https://www.betexplorer.com/football/portugal/liga-portugal-2/cf-os-belenenses-santa-clara/YDjeBj9C/"},

    @{r=131; idx=130; dataPartida=45283.57291666666; home="Maritimo"; homeGols=0; away="Penafiel"; awayGols=0;
      hOpen=1.61; hOpenDt="17/12/2023 15:12"; hClose=1.54; hCloseDt="23/12/2023 13:40";
      dOpen=3.83; dOpenDt="17/12/2023 15:12"; dClose=4.17; dCloseDt="23/12/2023 13:40";
      aOpen=5.29; aOpenDt="17/12/2023 15:12"; aClose=6.78; aCloseDt="23/12/2023 13:40";
      url="https://www.betexplorer.com/football/portugal/liga-portugal-2/maritimo-penafiel/bebrEhPg/"},

    @{r=132; idx=131; dataPartida=45290.5; home="FC Porto B"; homeGols=2; away="Nacional"; awayGols=3;
      hOpen=2.16; hOpenDt="23/12/2023 12:12"; hClose=2.07; hCloseDt="30/12/2023 11:58";
      dOpen=3.5; dOpenDt="23/12/2023 12:12"; dClose=3.89; dCloseDt="30/12/2023 11:58";
      aOpen=3.17; aOpenDt="23/12/2023 12:12"; aClose=3.41; aCloseDt="30/12/2023 11:58";
      url="https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-nacional/nZLADtjt/"},

    @{r=133; idx=132; dataPartida=45290.5; home="Leixoes"; homeGols=1; away="Academico Viseu"; awayGols=2;
      hOpen=2.82; hOpenDt="23/12/2023 12:12"; hClose=3.13; hCloseDt="30/12/2023 11:57";
      dOpen=3.15; dOpenDt="23/12/2023 12:12"; dClose=3.04; dCloseDt="30/12/2023 11:33";
      aOpen=2.58; aOpenDt="23/12/2023 12:12"; aClose=2.6; aCloseDt="30/12/2023 11:57";
      url="https://www.betexplorer.com/football/portugal/liga-portugal-2/leixoes-academico-viseu/fNiiCWg6/"},

    @{r=134; idx=133; dataPartida=45290.625; home="Pacos Ferreira"; homeGols=2; away="Benfica B"; awayGols=2;
      hOpen=2.02; hOpenDt="23/12/2023 15:12"; hClose=2.19; hCloseDt="30/12/2023 15:00";
      dOpen=3.41; dOpenDt="23/12/2023 15:12"; dClose=3.21; dCloseDt="30/12/2023 15:00";
      aOpen=3.59; aOpenDt="23/12/2023 15:12"; aClose=3.76; aCloseDt="30/12/2023 15:00";
      url="https://www.betexplorer.com/football/portugal/liga-portugal-2/pacos-ferreira-benfica/tdLEC05n/"},

    @{r=135; idx=134; dataPartida=45290.625; home="Tondela"; homeGols=1; away="Oliveirense"; awayGols=0;
      hOpen=1.5; hOpenDt="23/12/2023 12:12"; hClose=1.68; hCloseDt="30/12/2023 14:51";
      dOpen=4.22; dOpenDt="23/12/2023 12:12"; dClose=3.94; dCloseDt="30/12/2023 14:59";
      aOpen=6.04; aOpenDt="23/12/2023 12:12"; aClose=5.31; aCloseDt="30/12/2023 14:59";
      url="https://www.betexplorer.com/football/portugal/liga-portugal-2/tondela-oliveirense/ADUrIMsP/"},

    @{r=136; idx=135; dataPartida=45290.79166666666; home="Leiria"; homeGols=1; away="Feirense"; awayGols=1;
      hOpen=1.59; hOpenDt="23/12/2023 19:12"; hClose=1.76; hCloseDt="30/12/2023 18:52";
      dOpen=3.83; dOpenDt="23/12/2023 19:12"; dClose=3.57; dCloseDt="30/12/2023 18:59";
      aOpen=5.62; aOpenDt="23/12/2023 19:12"; aClose=5.28; aCloseDt="30/12/2023 18:59";
      url="https://www.betexplorer.com/football/portugal/liga-portugal-2/leiria-feirense/6u429UwP/"},

    @{r=137; idx=136; dataPartida=45291.5; home="Mafra"; homeGols=0; away="AVS"; awayGols=2;
      hOpen=2.51; hOpenDt="24/12/2023 12:11"; hClose=2.57; hCloseDt="31/12/2023 11:58";
      dOpen=3.13; dOpenDt="24/12/2023 12:11"; dClose=3.29; dCloseDt="31/12/2023 11:58";
      aOpen=2.89; aOpenDt="24/12/2023 12:11"; aClose=2.94; aCloseDt="31/12/2023 11:58";
      url="https://www.betexplorer.com/football/portugal/liga-portugal-2/mafra-avs/Sl5bAAOI/"}
)

foreach ($d in $newRows) {
    $r = $d.r
    $ws.Range("A$r").Value = $d.idx
    $ws.Range("B$r").Value = "portugal"
    $ws.Range("C$r").Value = "liga-portugal-2"
    $ws.Range("D$r").Value = "2023-2024"
    $ws.Range("E$r").Value = $d.dataPartida
    $ws.Range("F$r").Value = $d.home
    $ws.Range("G$r").Value = $d.homeGols
    $ws.Range("H$r").Value = $d.away
    $ws.Range("I$r").Value = $d.awayGols
    $ws.Range("J$r").Value = $d.hOpen
    $ws.Range("K$r").Value = $d.hOpenDt
    $ws.Range("L$r").Value = $d.hClose
    $ws.Range("M$r").Value = $d.hCloseDt
    $ws.Range("N$r").Value = $d.dOpen
    $ws.Range("O$r").Value = $d.dOpenDt
    $ws.Range("P$r").Value = $d.dClose
    $ws.Range("Q$r").Value = $d.dCloseDt
    $ws.Range("R$r").Value = $d.aOpen
    $ws.Range("S$r").Value = $d.aOpenDt
    $ws.Range("T$r").Value = $d.aClose
    $ws.Range("U$r").Value = $d.aCloseDt
    $ws.Range("V$r").Value = $d.url
}
